$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.138.28'
$ws.Range("E2").Value = '  +2.94%  '
$ws.Range("D3").Value = '3.066.27'
$ws.Range("E3").Value = '  +5.39%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '513.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.79%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.433'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.21'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("E10").Value = '  +3.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.371'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.14%  '
$ws.Range("D12").Value = '3.595.10'
$ws.Range("E12").Value = '  +5.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000164'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.20%  '
$ws.Range("D16").Value = '57.193.92'
$ws.Range("E16").Value = '  +3.16%  '
$ws.Range("D17").Value = '3.073.16'
$ws.Range("E17").Value = '  +5.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '335.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.98%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.499'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.41%  '
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").Value = '0.0₃0940'
$ws.Range("E27").Value = '  +11.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.81%  '
$ws.Range("E30").Value = '  +3.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '153.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0670'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.40%  '
$ws.Range("D39").Value = '3.104.01'
$ws.Range("E39").Value = '  +5.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.669'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.57%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.21%  '
$ws.Range("D44").Value = '2.238.35'
$ws.Range("E44").Value = '  +7.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0251'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.947'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0864'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.76%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.180'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.81%  '
